# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its text representation so values such as
# "1.000" or "27.272.38" are not reinterpreted by Excel as numbers/dates.
$ws.Range("D2:D51").NumberFormat = "@"

$priceUpdates = @{
    2 = "27.272.38"
    3 = "1.820.91"
    5 = "313.98"
    7 = "0.4479"
    8 = "0.3783"
    9 = "0.07426"
    10 = "0.8842"
    11 = "20.91"
    12 = "1.820.53"
    13 = "6.722"
    14 = "5.444"
    15 = "93.23"
    16 = "0.07130"
    18 = "0.000008805"
    19 = "1.000"
    20 = "15.11"
    21 = "27.299.60"
    22 = "5.372"
    23 = "10.92"
    24 = "1.960"
    25 = "151.69"
    26 = "2.303"
    27 = "18.66"
    28 = "5.355"
    29 = "117.46"
    30 = "0.08894"
    31 = "0.7877"
    32 = "1.198"
    33 = "4.593"
    34 = "2.911"
    36 = "1.112"
    37 = "0.01977"
    38 = "0.05280"
    39 = "7.340"
    40 = "0.5318"
    41 = "2.867"
    42 = "0.1708"
    43 = "2.295"
    44 = "8.615"
    45 = "0.5066"
    46 = "10.63"
    47 = "1.691"
    48 = "104.95"
    50 = "0.06397"
    51 = "66.05"
}

$volumeUpdates = @{
    2 = "  +0.87%  "
    3 = "  -0.05%  "
    4 = "  +0.01%  "
    5 = "  +0.69%  "
    6 = "  +0.00%  "
    7 = "  -1.82%  "
    8 = "  +1.92%  "
    9 = "  +1.92%  "
    10 = "  +3.11%  "
    11 = "  +0.49%  "
    12 = "  -0.01%  "
    13 = "  +1.02%  "
    14 = "  +2.11%  "
    15 = "  +0.53%  "
    16 = "  +0.51%  "
    17 = "  -0.07%  "
    18 = "  -0.14%  "
    19 = "  -0.01%  "
    20 = "  +0.82%  "
    21 = "  +1.02%  "
    22 = "  +3.85%  "
    23 = "  -0.29%  "
    24 = "  -1.51%  "
    25 = "  +0.03%  "
    26 = "  +3.92%  "
    27 = "  +1.23%  "
    28 = "  +1.98%  "
    29 = "  +0.59%  "
    30 = "  +0.07%  "
    31 = "  +4.64%  "
    32 = "  +0.65%  "
    33 = "  +3.03%  "
    34 = "  -0.96%  "
    35 = "  -0.03%  "
    36 = "  +1.33%  "
    37 = "  +0.59%  "
    38 = "  +0.59%  "
    39 = "  +2.00%  "
    40 = "  -0.01%  "
    41 = "  -0.35%  "
    43 = "  +16.74%  "
    44 = "  +0.91%  "
    45 = "  -2.77%  "
    46 = "  -0.20%  "
    47 = "  +1.16%  "
    48 = "  -0.52%  "
    49 = "  +0.00%  "
    50 = "  +0.20%  "
    51 = "  +4.32%  "
}

foreach ($row in $priceUpdates.Keys) {
    $ws.Cells.Item($row, 4).Value = $priceUpdates[$row]
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Cells.Item($row, 5).Value = $volumeUpdates[$row]
}
